$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3497718.5
$ws.Range("I33").Value = 1639.7333
$ws.Range("J33").Value = 8265099
$ws.Range("K33").Value = 1639.7333
$ws.Range("L33").Value = 8265099
$ws.Range("M33").Value = -1410.7333
$ws.Range("N33").Value = -8265557

$ws.Range("H86").Value = 9966.666999999999
$ws.Range("I86").Value = 9966.666999999999
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 9966.666999999999
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -8843.666999999999
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 9966.666999999999
$ws.Range("I89").Value = 9966.666999999999
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 49833.335
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -44217.335
$ws.Range("N89").ClearContents()

$ws.Range("H107").Value = 35715212
$ws.Range("I107").Value = 35715212
$ws.Range("K107").Value = 35715212
$ws.Range("M107").Value = -35713292

$ws.Range("H129").Value = 1402.5
$ws.Range("J129").Value = 1847.125
$ws.Range("L129").Value = 5541.375
$ws.Range("N129").Value = -15541.375

$ws.Range("H136").Value = 45797
$ws.Range("J136").Value = 45797
$ws.Range("L136").Value = 45797
$ws.Range("N136").Value = -55997

$ws.Range("H138").Value = 2473.141
$ws.Range("I138").Value = 1166.3043
$ws.Range("J138").Value = 4351.7188
$ws.Range("K138").Value = 3498.9129
$ws.Range("L138").Value = 13055.1564
$ws.Range("M138").Value = 1641.0871
$ws.Range("N138").Value = -23335.1564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10651.9
$ws.Range("I32").Value = 7969.4805
$ws.Range("J32").Value = 19632.174
$ws.Range("K32").Value = 7969.4805
$ws.Range("L32").Value = 19632.174
$ws.Range("M32").Value = -7682.4805
$ws.Range("N32").Value = -20206.174

$ws.Range("H74").Value = 7577024
$ws.Range("I74").Value = 1035.7435
$ws.Range("K74").Value = 1035.7435
$ws.Range("M74").Value = -161.7435

$ws.Range("H77").Value = 7577024
$ws.Range("I77").Value = 1035.7435
$ws.Range("K77").Value = 5178.717500000001
$ws.Range("M77").Value = -810.7175000000007

$ws.Range("H110").Value = 22562.4
$ws.Range("I110").Value = 24009.715
$ws.Range("J110").Value = 2300
$ws.Range("K110").Value = 24009.715
$ws.Range("L110").Value = 2300
$ws.Range("M110").Value = -21964.715
$ws.Range("N110").Value = -6390

$ws.Range("H133").Value = 28493.334
$ws.Range("J133").Value = 28493.334
$ws.Range("L133").Value = 28493.334
$ws.Range("N133").Value = -33553.334

$ws.Range("H138").Value = 44214.5
$ws.Range("J138").Value = 44214.5
$ws.Range("L138").Value = 44214.5
$ws.Range("N138").Value = -54494.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 166669010
$ws.Range("I99").Value = 250001120
$ws.Range("J99").Value = 4750
$ws.Range("K99").Value = 250001120
$ws.Range("L99").Value = 4750
$ws.Range("M99").Value = -249999622
$ws.Range("N99").Value = -7746

$ws.Range("H137").Value = 67668
$ws.Range("J137").Value = 57085
$ws.Range("L137").Value = 57085
$ws.Range("N137").Value = -67285

$ws.Range("H138").Value = 60437.145
$ws.Range("J138").Value = 60437.145
$ws.Range("L138").Value = 60437.145
$ws.Range("N138").Value = -70717.14499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1286.5

$ws.Range("H32").Value = 1591.25
$ws.Range("I32").Value = 1591.25
$ws.Range("K32").Value = 1591.25
$ws.Range("M32").Value = -1275.25

$ws.Range("H53").Value = 38500
$ws.Range("J53").Value = 38500
$ws.Range("L53").Value = 38500
$ws.Range("N53").Value = -39714

$ws.Range("H99").Value = 25001020
$ws.Range("I99").Value = 1275
$ws.Range("K99").Value = 1275
$ws.Range("M99").Value = 223

$ws.Range("H107").Value = 1027.7059
$ws.Range("I107").Value = 683.9231
$ws.Range("J107").Value = 2145
$ws.Range("K107").Value = 683.9231
$ws.Range("L107").Value = 2145
$ws.Range("M107").Value = 1236.0769
$ws.Range("N107").Value = -5985

$ws.Range("H126").Value = 25001020
$ws.Range("I126").Value = 1275
$ws.Range("K126").Value = 3825
$ws.Range("M126").Value = -1355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1977849
$ws.Range("I5").Value = 438.81818
$ws.Range("J5").Value = 3790475
$ws.Range("K5").Value = 1316.45454
$ws.Range("L5").Value = 11371425
$ws.Range("M5").Value = -1204.45454
$ws.Range("N5").Value = -11371649

$ws.Range("H50").Value = 41.25
$ws.Range("I50").Value = 35.714287
$ws.Range("K50").Value = 107.142861
$ws.Range("M50").Value = 373.857139

$ws.Range("H53").Value = 41.25
$ws.Range("I53").Value = 35.714287
$ws.Range("K53").Value = 107.142861
$ws.Range("M53").Value = 373.857139

$ws.Range("H132").Value = 7264470
$ws.Range("I132").Value = 1951
$ws.Range("J132").Value = 9339476
$ws.Range("K132").Value = 17559
$ws.Range("L132").Value = 84055284
$ws.Range("M132").Value = -15029
$ws.Range("N132").Value = -84060344

$ws.Range("H135").Value = 1977849
$ws.Range("I135").Value = 438.81818
$ws.Range("J135").Value = 3790475
$ws.Range("K135").Value = 3949.36362
$ws.Range("L135").Value = 34114275
$ws.Range("M135").Value = -1414.36362
$ws.Range("N135").Value = -34119345

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H97").Value = 1910
$ws.Range("I97").Value = 1910
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1910
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1414
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 2018.174
$ws.Range("I102").Value = 1618.125
$ws.Range("J102").Value = 2932.5715
$ws.Range("K102").Value = 1618.125
$ws.Range("L102").Value = 2932.5715
$ws.Range("M102").Value = 3.875
$ws.Range("N102").Value = -6176.5715

$ws.Range("H138").Value = 39491.5
$ws.Range("J138").Value = 39491.5
$ws.Range("L138").Value = 39491.5
$ws.Range("N138").Value = -49771.5

$ws.Range("H141").Value = 57889.43
$ws.Range("J141").Value = 57889.43
$ws.Range("L141").Value = 57889.43
$ws.Range("N141").Value = -68249.42999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2761.2727
$ws.Range("I61").Value = 2396.75
$ws.Range("K61").Value = 2396.75
$ws.Range("M61").Value = -2194.75

$ws.Range("H93").Value = 1428.8572
$ws.Range("I93").Value = 1000.4
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1000.4
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = 247.6
$ws.Range("N93").Value = -4996

$ws.Range("H113").Value = 2761.2727
$ws.Range("I113").Value = 2396.75
$ws.Range("K113").Value = 2396.75
$ws.Range("M113").Value = -226.75

$ws.Range("H122").Value = 7409772.5
$ws.Range("I122").Value = 14297359
$ws.Range("K122").Value = 42892077
$ws.Range("M122").Value = -42889627

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 37494
$ws.Range("J46").Value = 37494
$ws.Range("L46").Value = 37494
$ws.Range("N46").Value = -37956

$ws.Range("H132").Value = 1952.4062
$ws.Range("I132").Value = 877.4
$ws.Range("J132").Value = 2441.0454
$ws.Range("K132").Value = 2632.2
$ws.Range("L132").Value = 7323.1362
$ws.Range("M132").Value = -102.1999999999998
$ws.Range("N132").Value = -12383.1362

$ws.Range("H134").Value = 37494
$ws.Range("J134").Value = 37494
$ws.Range("L134").Value = 112482
$ws.Range("N134").Value = -117552

$ws.Range("H136").Value = 4277066.5
$ws.Range("I136").Value = 4377.5293
$ws.Range("J136").Value = 7578690
$ws.Range("K136").Value = 13132.5879
$ws.Range("L136").Value = 22736070
$ws.Range("M136").Value = -10582.5879
$ws.Range("N136").Value = -22741170
